$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "data_voo" metadata row (original row 2) entirely; everything
# below it shifts up by one row.
$ws.Rows(2).Delete()

# The "arr_time" row (now row 3) had its raw_null_tolerance value corrected
# from 0.1 to 0.05.
$ws.Range("G3").Value = 0.05

# Move the selection to reflect where the user ended up after editing.
$ws.Range("F10").Select()
